# Applies per-cell value updates for the crypto price/volume table.
# Cells whose new text would otherwise be auto-parsed by Excel as a
# number (losing formatting like trailing zeros, e.g. "0.999" or "1.40")
# are briefly switched to Text format, written, then restored to their
# original style so the cell keeps behaving exactly as before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.070.68"
$ws.Range("E2").Value = "  -2.84%  "
$ws.Range("D3").Value = "2.372.73"
$ws.Range("E3").Value = "  -4.97%  "
$ws.Range("E4").Value = "  +0.06%  "
$sty = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.08"
$ws.Range("D5").Style = $sty
$ws.Range("E5").Value = "  -2.14%  "
$sty = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.33"
$ws.Range("D6").Style = $sty
$ws.Range("E6").Value = "  -4.94%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -5.52%  "
$ws.Range("D9").Value = "2.373.45"
$ws.Range("E9").Value = "  -4.93%  "
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("E13").Value = "  -4.56%  "
$sty = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.15"
$ws.Range("D14").Style = $sty
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "2.809.69"
$ws.Range("E15").Value = "  -4.68%  "
$ws.Range("D17").Value = "59.998.74"
$ws.Range("D18").Value = "2.383.73"
$ws.Range("E18").Value = "  -4.45%  "
$sty = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.54"
$ws.Range("D19").Style = $sty
$ws.Range("E19").Value = "  -4.98%  "
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$sty = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.61"
$ws.Range("D21").Style = $sty
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$sty = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.65"
$ws.Range("D22").Style = $sty
$ws.Range("E22").Value = "  -5.18%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +3.21%  "
$sty = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.78"
$ws.Range("D25").Style = $sty
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "2.503.72"
$ws.Range("E27").Value = "  -4.52%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0914"
$ws.Range("E28").Value = "  -8.63%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$sty = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.64"
$ws.Range("D29").Style = $sty
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$sty = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.97"
$ws.Range("D30").Style = $sty
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$sty = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.40"
$ws.Range("D31").Style = $sty
$ws.Range("E31").Value = "  -5.28%  "
$sty = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "506.84"
$ws.Range("D32").Style = $sty
$ws.Range("E32").Value = "  -5.03%  "
$ws.Range("E33").Value = "  -4.73%  "
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("E35").Value = "  -0.84%  "
$sty = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = $sty
$ws.Range("E36").Value = "  -0.05%  "
$sty = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.60"
$ws.Range("D37").Style = $sty
$ws.Range("E37").Value = "  -5.59%  "
$sty = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.38"
$ws.Range("D38").Style = $sty
$ws.Range("E38").Value = "  -8.41%  "
$ws.Range("E39").Value = "  -1.89%  "
$sty = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.92"
$ws.Range("D40").Style = $sty
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  +0.40%  "
$sty = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.05"
$ws.Range("D43").Style = $sty
$ws.Range("E43").Value = "  -4.85%  "
$sty = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.23"
$ws.Range("D44").Style = $sty
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("E45").Value = "  -5.39%  "
$sty = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "139.45"
$ws.Range("D46").Style = $sty
$ws.Range("E46").Value = "  -6.20%  "
$ws.Range("E47").Value = "  -1.85%  "
$sty = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.12"
$ws.Range("D48").Style = $sty
$ws.Range("E49").Value = "  -3.94%  "
$sty = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.573"
$ws.Range("D50").Style = $sty
$ws.Range("E50").Value = "  -2.53%  "
$sty = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0920"
$ws.Range("D51").Style = $sty
$ws.Range("E51").Value = "  -2.30%  "
